$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(75, 'Thu Oct 17 00:37:19 2024', 150, 144),
    @(76, 'Thu Oct 17 00:37:29 2024', 147, 145),
    @(77, 'Thu Oct 17 00:46:51 2024', 131, 136),
    @(78, 'Thu Oct 17 00:46:51 2024', 132, 137),
    @(79, 'Thu Oct 17 00:46:51 2024', 133, 138),
    @(80, 'Thu Oct 17 00:46:55 2024', 131, 137),
    @(81, 'Thu Oct 17 00:47:05 2024', 131, 138),
    @(82, 'Thu Oct 17 00:47:15 2024', 130, 137),
    @(83, 'Thu Oct 17 00:50:38 2024', 128, 134),
    @(84, 'Thu Oct 17 00:50:48 2024', 128, 133),
    @(85, 'Thu Oct 17 01:07:40 2024', 124, 132),
    @(86, 'Thu Oct 17 01:07:50 2024', 123, 132),
    @(87, 'Thu Oct 17 01:08:00 2024', 124, 133),
    @(88, 'Thu Oct 17 01:14:41 2024', 124, 133),
    @(89, 'Thu Oct 17 01:14:51 2024', 124, 133),
    @(90, 'Thu Oct 17 01:15:01 2024', 125, 134),
    @(91, 'Thu Oct 17 01:15:11 2024', 297, 134),
    @(92, 'Thu Oct 17 01:16:35 2024', 138, 132),
    @(93, 'Thu Oct 17 01:16:45 2024', 135, 131),
    @(94, 'Thu Oct 17 01:16:55 2024', 135, 133),
    @(95, 'Thu Oct 17 01:17:05 2024', 134, 134),
    @(96, 'Thu Oct 17 01:17:15 2024', 132, 135),
    @(97, 'Thu Oct 17 01:17:25 2024', 134, 135),
    @(98, 'Thu Oct 17 01:17:35 2024', 131, 133),
    @(99, 'Thu Oct 17 01:17:45 2024', 131, 134),
    @(100, 'Thu Oct 17 01:17:55 2024', 130, 133),
    @(101, 'Thu Oct 17 01:18:05 2024', 133, 137),
    @(102, 'Thu Oct 17 01:18:15 2024', 130, 134),
    @(103, 'Thu Oct 17 01:18:25 2024', 128, 133),
    @(104, 'Thu Oct 17 01:18:35 2024', 130, 134),
    @(105, 'Thu Oct 17 01:18:45 2024', 130, 136),
    @(106, 'Thu Oct 17 01:18:55 2024', 127, 132),
    @(107, 'Thu Oct 17 01:19:05 2024', 126, 132),
    @(108, 'Thu Oct 17 01:19:15 2024', 125, 131),
    @(109, 'Thu Oct 17 01:19:25 2024', 127, 133),
    @(110, 'Thu Oct 17 01:19:35 2024', 128, 135),
    @(111, 'Thu Oct 17 01:19:45 2024', 128, 135),
    @(112, 'Thu Oct 17 01:19:55 2024', 127, 133),
    @(113, 'Thu Oct 17 01:20:05 2024', 129, 135),
    @(114, 'Thu Oct 17 01:20:15 2024', 130, 136),
    @(115, 'Thu Oct 17 01:20:25 2024', 127, 134),
    @(116, 'Thu Oct 17 01:20:35 2024', 127, 134),
    @(117, 'Thu Oct 17 01:20:45 2024', 130, 137),
    @(118, 'Thu Oct 17 01:21:00 2024', 127, 132),
    @(119, 'Thu Oct 17 01:21:10 2024', 162, 132),
    @(120, 'Thu Oct 17 01:21:20 2024', 187, 132),
    @(121, 'Thu Oct 17 01:21:30 2024', 162, 134),
    @(122, 'Thu Oct 17 01:21:40 2024', 151, 133),
    @(123, 'Thu Oct 17 01:21:50 2024', 145, 134),
    @(124, 'Thu Oct 17 01:22:00 2024', 141, 135),
    @(125, 'Thu Oct 17 01:22:10 2024', 136, 133),
    @(126, 'Thu Oct 17 01:22:20 2024', 135, 133),
    @(127, 'Thu Oct 17 01:22:30 2024', 134, 134),
    @(128, 'Thu Oct 17 01:22:40 2024', 133, 134),
    @(129, 'Thu Oct 17 01:22:50 2024', 132, 132),
    @(130, 'Thu Oct 17 01:23:00 2024', 131, 133),
    @(131, 'Thu Oct 17 01:23:10 2024', 128, 131),
    @(132, 'Thu Oct 17 01:23:20 2024', 126, 131),
    @(133, 'Thu Oct 17 01:23:30 2024', 125, 130),
    @(134, 'Thu Oct 17 01:23:40 2024', 124, 130)
)

$startRow = 76
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
